$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("A4").Value = "'2018.07.04"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "15:02:44"
$ws.Range("C4").Value = "RS"
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 0.1
$ws.Range("J4").Value = 0.0095
$ws.Range("K4").Value = "effective"
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 3.19
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = "N/A"
$ws.Range("P4").Value = "N/A"

# --- Row 5 ---
$ws.Range("A5").Value = "'2018.07.04"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "15:04:31"
$ws.Range("C5").Value = "RS"
$ws.Range("D5").Value = 32
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 0.1
$ws.Range("J5").Value = 0.96
$ws.Range("K5").Value = "effective"
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 8.56
$ws.Range("N5").Value = 111
$ws.Range("O5").Value = 36.7
$ws.Range("P5").Value = 0.5385884953938632

# --- Row 6 ---
$ws.Range("A6").Value = "'2018.07.04"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "15:04:50"
$ws.Range("C6").Value = "RS"
$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 4000
$ws.Range("G6").Value = 99
$ws.Range("H6").Value = 250
$ws.Range("I6").Value = 0.1
$ws.Range("J6").Value = 0.96
$ws.Range("K6").Value = "effective"
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 8.1
$ws.Range("N6").Value = 81
$ws.Range("O6").Value = 43.85
$ws.Range("P6").Value = 0.4066105357283207
